$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 254, shifting existing rows 254-338 down to 255-339
$ws.Rows(254).Insert()

# Populate the newly inserted row 254 with the new record's data
$ws.Range("A254").Value = 9
$ws.Range("B254").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C254").Value = "Metropolitana"
$ws.Range("D254").Value = 44876
$ws.Range("E254").Value = 13
$ws.Range("F254").Value = 300000001
$ws.Range("G254").Value = "Rabanito"
$ws.Range("H254").Value = "Sin especificar"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 12000
$ws.Range("K254").Value = 3000
$ws.Range("L254").Value = 4000
$ws.Range("M254").Value = 3417
$ws.Range("N254").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O254").Value = "Provincia de Chacabuco"
$ws.Range("P254").Value = 34
$ws.Range("Q254").Value = 100
$ws.Range("R254").Value = "Hortaliza"
